$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "2025-04-28 21:16:18"
$ws.Range("B19").Value = 13
